$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) contains numeric-looking text values (e.g. "1.00", "8.59").
# Excel COM auto-converts such strings to real numbers on assignment, which would
# strip meaningful trailing zeros / formatting from the original text cells.
# Force the whole column to Text format first so the new values stay text exactly
# as scraped, then restore the default (unstyled) look so no visible formatting changes.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.771.29"
$ws.Range("E2").Value = "  -1.01%  "

$ws.Range("D3").Value = "2.365.36"
$ws.Range("E3").Value = "  -1.18%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "317.37"
$ws.Range("E5").Value = "  -2.96%  "

$ws.Range("D6").Value = "108.89"
$ws.Range("E6").Value = "  +3.23%  "

$ws.Range("D7").Value = "0.635"
$ws.Range("E7").Value = "  -2.08%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  -4.53%  "

$ws.Range("D10").Value = "41.96"
$ws.Range("E10").Value = "  -0.03%  "

$ws.Range("D11").Value = "0.0925"
$ws.Range("E11").Value = "  -1.51%  "

$ws.Range("D12").Value = "8.59"
$ws.Range("E12").Value = "  -1.33%  "

$ws.Range("D13").Value = "1.00"
$ws.Range("E13").Value = "  -5.52%  "

$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("D15").Value = "16.12"
$ws.Range("E15").Value = "  -5.47%  "

$ws.Range("D16").Value = "2.723.73"
$ws.Range("E16").Value = "  -1.01%  "

$ws.Range("D17").Value = "2.379.21"
$ws.Range("E17").Value = "  -0.24%  "

$ws.Range("D18").Value = "42.739.62"
$ws.Range("E18").Value = "  -1.03%  "

$ws.Range("D19").Value = "7.73"
$ws.Range("E19").Value = "  -0.19%  "

$ws.Range("E20").Value = "  -1.74%  "

$ws.Range("D21").Value = "76.25"
$ws.Range("E21").Value = "  -0.68%  "

$ws.Range("D22").Value = "3.65"
$ws.Range("E22").Value = "  -3.27%  "

$ws.Range("D23").Value = "256.16"
$ws.Range("E23").Value = "  -6.65%  "

$ws.Range("D24").Value = "2.32"
$ws.Range("E24").Value = "  -4.53%  "

$ws.Range("D25").Value = "9.52"
$ws.Range("E25").Value = "  -0.53%  "

$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").Value = "11.51"
$ws.Range("E27").Value = "  -2.43%  "

$ws.Range("D28").Value = "22.89"
$ws.Range("E28").Value = "  -1.20%  "

$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  +2.52%  "

$ws.Range("D30").Value = "37.41"
$ws.Range("E30").Value = "  +0.33%  "

$ws.Range("D31").Value = "171.46"
$ws.Range("E31").Value = "  -2.57%  "

$ws.Range("E32").Value = "  -4.46%  "

$ws.Range("D33").Value = "6.08"
$ws.Range("E33").Value = "  +2.10%  "

$ws.Range("E34").Value = "  -8.37%  "

$ws.Range("D35").Value = "0.123"
$ws.Range("E35").Value = "  +13.08%  "

$ws.Range("E36").Value = "  -2.36%  "

$ws.Range("D37").Value = "4.70"
$ws.Range("E37").Value = "  -3.97%  "

$ws.Range("D38").Value = "0.0364"
$ws.Range("E38").Value = "  -0.28%  "

$ws.Range("D39").Value = "3.91"
$ws.Range("E39").Value = "  -5.81%  "

$ws.Range("D40").Value = "2.69"
$ws.Range("E40").Value = "  -4.98%  "

$ws.Range("D41").Value = "0.243"
$ws.Range("E41").Value = "  +4.06%  "

$ws.Range("D42").Value = "1.51"
$ws.Range("E42").Value = "  -5.15%  "

$ws.Range("D43").Value = "71.22"
$ws.Range("E43").Value = "  +1.62%  "

$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("D45").Value = "12.35"
$ws.Range("E45").Value = "  +0.19%  "

$ws.Range("B46").Value = "THORChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D46").Value = "5.58"
$ws.Range("E46").Value = "  +0.62%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "112.40"
$ws.Range("E47").Value = "  -8.01%  "

$ws.Range("D48").Value = "9.24"
$ws.Range("E48").Value = "  -1.35%  "

$ws.Range("D49").Value = "86.09"
$ws.Range("E49").Value = "  -3.83%  "

$ws.Range("D50").Value = "77.23"
$ws.Range("E50").Value = "  +5.35%  "

$ws.Range("E51").Value = "  -0.92%  "

$ws.Range("D2:D51").Style = "Normal"